# Backfill profile.equilibrated column (and more)
# - Replace boolean borehole.equilibrated (measurement sheet column H) with
#   categorical borehole.equilibrium: true, estimated, false (see #96)

$wb = $excel.ActiveWorkbook

$meas  = $wb.Worksheets.Item("measurement")
$lists = $wb.Worksheets.Item("lists")

# 1. Rename the "measurement" sheet header H1 from equilibrated -> equilibrium
$meas.Range("H1").Value = "equilibrium"

# 2. Update the H1 cell comment to describe the new categorical field
$newComment = "[string] equilibrium`nWhether and how reported temperatures equilibrated following drilling.`n- true: Equilibrium was measured`n- estimated: Equilibrium was estimated (typically by extrapolation)`n- false: Equilibrium was not reached`nconstraints:`n  - enum: ['true', 'estimated', 'false']"
$meas.Range("H1").Comment.Text($newComment)

# 3. Extend the hidden "lists" sheet with a third column of allowed
#    equilibrium values (true / estimated / false). The leading apostrophe
#    forces "true"/"false" to be stored as text instead of being coerced to
#    native booleans; resetting the style back to Normal afterwards drops
#    the quote-prefix marker so the cell matches a plain text cell.
$lists.Range("C1").Value = "'true"
$lists.Range("C1").Style = "Normal"
$lists.Range("C2").Value = "estimated"
$lists.Range("C3").Value = "'false"
$lists.Range("C3").Style = "Normal"

# 4. Point the H column conditional formatting at the new list instead of
#    the old hard-coded TRUE/FALSE comparison
$rngH = $meas.Range("H2:H1048576")
$cond = $rngH.FormatConditions.Item(1)
$cond.Formula1 = "=IF(ISBLANK(H2), FALSE, ISNA(MATCH(H2, 'lists'!`$C`$1:`$C`$3, 0)))"

# 5. Point the H column data validation dropdown at the new list
$dv = $rngH.Validation
$dv.Formula1 = "='lists'!`$C`$1:`$C`$3"

# 6. Match the H column width to the new (narrower) label width
$meas.Columns.Item(8).ColumnWidth = 13.0

Write-Host "equilibrium backfill applied"
